$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")
$ws.Activate()

# Insert a new row at row 8, pushing the old row 8 ("Complete project plan...") down to row 9.
# Excel copies formatting from the row above (row 7) onto the new row.
$ws.Rows("8").Insert()

# Fill in the new Row 8 with the "Complete RTM" milestone data
$ws.Range("A8").Value = "Complete RTM"
$ws.Range("B8").Value = 43867
$ws.Range("C8").Value = 43867
$ws.Range("D8").Value = 43867
$ws.Range("F8").Value = "Currently working on modifications asked by Eng. Mohamed Ali"

# E8 stays empty but gets a plain wrap-text alignment (no horizontal/vertical centering)
$ws.Range("E8").ClearFormats()
$ws.Range("E8").WrapText = $true

# Fill in Actual End Date (E5, E6) with the finished dates, reusing the
# existing date-formatted style from the Planned End Date column (D)
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 43866

$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 43866

# Resize the table to include the new row
$ws.ListObjects("Table2").Resize($ws.Range("A4:F9"))

# Update the selection to match the target state
$ws.Range("A8").Select()
